$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TODO Before 0.0.1")
$ws2 = $wb.Worksheets.Item("TODO's")
$ws3 = $wb.Worksheets.Item("Logs")

# --- Sheet "Logs": rework the stone-throwing log entry to mention crystals too ---
$ws3.Range("B69").Value = "rework stone and crystal throwing to rigid body "

# --- Sheet "TODO Before 0.0.1": add two new TODO rows (47, 48) ---
$ws1.Range("A47").Value = 44
$ws1.Range("B47").Value = "add map to folow player oppened place"
$ws1.Range("C47").Value = "in-progress"
$ws1.Range("C47").HorizontalAlignment = -4108
$ws1.Range("C47").VerticalAlignment = -4108

$ws1.Range("A48").Value = 45
$ws1.Range("B48").Value = "think to get rid of gold coins concept - better sould to be used"
$ws1.Range("C48").Value = "in-progress"
$ws1.Range("C48").HorizontalAlignment = -4108
$ws1.Range("C48").VerticalAlignment = -4108

# --- Sheet "TODO's": widen column D ---
$ws2.Columns.Item(4).ColumnWidth = 14.166666666666666

# --- Selections (apply on non-active sheets first, then finish on the active sheet) ---
$ws2.Range("D6").Select()
$ws3.Range("B70").Select()

$ws1.Activate()
$ws1.Range("C48").Select()
